# Correct the ABUNDANCE column values (column 2) of Table 7 per the
# "corrected abundance sum error" fix.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    2  = "0.600"   # CODNEAR
    3  = "0.339"   # CODNEARNCW
    4  = "0.348"   # CODFAPL
    5  = "0.583"   # CODICE
    6  = "0.199"   # CODBA2532
    7  = "0.408"   # CODKAT
    8  = "0.092"   # CODIS
    9  = "0.204"   # CODVIa
    10 = "0.329"   # CODIIIaW
    11 = "0.764"   # HAKENRTN
    12 = "0.427"   # HAKESOTH
}

foreach ($rowIndex in $newValues.Keys) {
    $cell = $t.Cell($rowIndex, 2)
    $cellRange = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, preserving the cell's run formatting.
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $newValues[$rowIndex]
}
